$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Put the new text value into A1 (creates/uses the shared-string table entry).
$ws.Range("A1").Value = "aaaaaaaaaa"

# Switch the default/normal font from Arial to Calibri.
$ws.Cells.Font.Name = "Calibri"
$excel.StandardFont = "Calibri"
